$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.289, 0.046, 0.214, 0.168, 0.028, 0.168, 0.2, 0.026, 0.161, 0.307, 0.099, 0.314, 0.254, 0.022, 0.149, 0.478, 0.222, 0.471, 0.244, 0.08400000000000001, 0.289, 0.242, 0.043, 0.206, 0.431, 0.13, 0.361, 0.117, 0.007, 0.082, 0.713, 0.106, 0.325, 0.656, 0.165, 0.406, 0.671, 0.121, 0.348, 0.68),
    @(0.822, 0.146, 0.382, 0.667, 0.222, 0.471, 0.8, 0.16, 0.4, 0.6, 0.24, 0.49, 0.8, 0.16, 0.4, 0.533, 0.249, 0.499, 0.511, 0.25, 0.5, 0.733, 0.196, 0.442, 0.8, 0.16, 0.4, 0.711, 0.205, 0.453, 0.956, 0.042, 0.206, 0.778, 0.173, 0.416, 0.911, 0.081, 0.285, 0.882),
    @(0.428, $null, $null, 0.268, $null, $null, 0.32, $null, $null, 0.406, $null, $null, 0.386, $null, $null, 0.504, $null, $null, 0.33, $null, $null, 0.364, $null, $null, 0.5600000000000001, $null, $null, 0.201, $null, $null, 0.8169999999999999, $null, $null, 0.712, $null, $null, 0.773, $null, $null, 0.767),
    @(0.601, $null, $null, 0.418, $null, $null, 0.5, $null, $null, 0.504, $null, $null, 0.5590000000000001, $null, $null, 0.521, $null, $null, 0.419, $null, $null, 0.521, $null, $null, 0.6830000000000001, $null, $null, 0.353, $null, $null, 0.895, $null, $null, 0.75, $null, $null, 0.85, $null, $null, 0.832),
    @(0.742, 0.149, 0.386, 0.5629999999999999, 0.192, 0.439, 0.697, 0.158, 0.398, 0.531, 0.211, 0.46, 0.713, 0.157, 0.396, 0.509, 0.234, 0.484, 0.445, 0.21, 0.458, 0.662, 0.182, 0.426, 0.737, 0.157, 0.396, 0.596, 0.189, 0.435, 0.879, 0.06, 0.244, 0.77, 0.172, 0.415, 0.878, 0.08599999999999999, 0.294, 0.842),
    @(0.644, 0.229, 0.479, 0.444, 0.247, 0.497, 0.578, 0.244, 0.494, 0.444, 0.247, 0.497, 0.6, 0.24, 0.49, 0.467, 0.249, 0.499, 0.356, 0.229, 0.479, 0.556, 0.247, 0.497, 0.644, 0.229, 0.479, 0.489, 0.25, 0.5, 0.756, 0.185, 0.43, 0.756, 0.185, 0.43, 0.822, 0.146, 0.382, 0.778),
    @(0.778, 0.173, 0.416, 0.6, 0.24, 0.49, 0.733, 0.196, 0.442, 0.6, 0.24, 0.49, 0.778, 0.173, 0.416, 0.533, 0.249, 0.499, 0.511, 0.25, 0.5, 0.733, 0.196, 0.442, 0.8, 0.16, 0.4, 0.6, 0.24, 0.49, 0.956, 0.042, 0.206, 0.778, 0.173, 0.416, 0.911, 0.081, 0.285, 0.882),
    @(0.822, 0.146, 0.382, 0.667, 0.222, 0.471, 0.8, 0.16, 0.4, 0.6, 0.24, 0.49, 0.8, 0.16, 0.4, 0.533, 0.249, 0.499, 0.511, 0.25, 0.5, 0.733, 0.196, 0.442, 0.8, 0.16, 0.4, 0.644, 0.229, 0.479, 0.956, 0.042, 0.206, 0.778, 0.173, 0.416, 0.911, 0.081, 0.285, 0.882),
    @(1.378, 0.668, 0.8169999999999999, 1.633, 1.032, 1.016, 1.556, 1.191, 1.091, 1.407, 0.538, 0.733, 1.389, 0.571, 0.756, 1.125, 0.109, 0.331, 1.435, 0.507, 0.712, 1.303, 0.332, 0.577, 1.25, 0.299, 0.546, 2, 3.812, 1.953, 1.233, 0.225, 0.474, 1.029, 0.028, 0.167, 1.098, 0.08799999999999999, 0.297, 1.12),
    @(3.533, 1.404, 1.185, 4.564, 0.707, 0.841, 4.524, 0.916, 0.957, 2.3, 0.61, 0.781, 3.333, 0.756, 0.869, 1.25, 0.188, 0.433, 2.517, 1.008, 1.004, 4.073, 6.751, 2.598, 2.833, 3.901, 1.975, 6.273, 2.88, 1.697, 1.667, 0.8, 0.894, 1.311, 0.348, 0.59, 1.689, 0.792, 0.89, 1.556)
)

$startRow = 4
$startCol = 2  # column B

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $val
        }
    }
}
